$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated NATMI LR-pair statistics (Grn-Tnfrsf1b), per commit "Natmi following Dr Hou advice".
# Ligand-/receptor-expressing cell counts moved from 1 to 3 cells, which also shifts the
# derived average/total expression, specificity and edge-weight columns for every data row.
$rowData = @{
    2 = @{ "E" = 3; "G" = 24.50527833333333; "H" = 73.515835; "I" = 0.04846830138877924; "J" = 0.04846830138877924; "K" = 3; "M" = 11.84297166666667; "N" = 35.528915; "O" = 0.1201507426883529; "P" = 0.1201507426883529; "Q" = 290.2153169854472; "R" = 2611.937852869024; "S" = 0.005823502408704753; "T" = 0.005823502408704754 }
    3 = @{ "E" = 3; "G" = 24.50527833333333; "H" = 73.515835; "I" = 0.04846830138877924; "J" = 0.04846830138877924; "K" = 3; "M" = 10.940628; "N" = 32.821884; "O" = 0.1109961770302011; "P" = 0.1109961770302011; "Q" = 268.10313428146; "R" = 2412.92820853314; "S" = 0.005379796161302082; "T" = 0.005379796161302083 }
    4 = @{ "E" = 3; "G" = 24.50527833333333; "H" = 73.515835; "I" = 0.04846830138877924; "J" = 0.04846830138877924; "K" = 3; "M" = 73.27090866666667; "N" = 219.812726; "O" = 0.7433568483938062; "P" = 0.7433568483938062; "Q" = 1795.524010612912; "R" = 16159.71609551621; "S" = 0.03602924376736408; "T" = 0.03602924376736408 }
    5 = @{ "E" = 3; "G" = 24.50527833333333; "H" = 73.515835; "I" = 0.04846830138877924; "J" = 0.04846830138877924; "K" = 3; "M" = 2.513102666666667; "N" = 7.539308; "O" = 0.02549623188763971; "P" = 0.02549623188763971; "Q" = 61.58428032690889; "R" = 554.25852294218; "S" = 0.001235759051408325; "T" = 0.001235759051408325 }
    6 = @{ "E" = 3; "G" = 110.1980973333333; "H" = 330.594292; "I" = 0.2179577200213544; "J" = 0.2179577200213544; "K" = 3; "M" = 11.84297166666667; "N" = 35.528915; "O" = 0.1201507426883529; "P" = 0.1201507426883529; "Q" = 1305.072944439242; "R" = 11745.65649995318; "S" = 0.02618778193522582; "T" = 0.02618778193522583 }
    7 = @{ "E" = 3; "G" = 110.1980973333333; "H" = 330.594292; "I" = 0.2179577200213544; "J" = 0.2179577200213544; "K" = 3; "M" = 10.940628; "N" = 32.821884; "O" = 0.1109961770302011; "P" = 0.1109961770302011; "Q" = 1205.636389231792; "R" = 10850.72750308613; "S" = 0.02419247367658927; "T" = 0.02419247367658927 }
    8 = @{ "E" = 3; "G" = 110.1980973333333; "H" = 330.594292; "I" = 0.2179577200213544; "J" = 0.2179577200213544; "K" = 3; "M" = 73.27090866666667; "N" = 219.812726; "O" = 0.7433568483938062; "P" = 0.7433568483938062; "Q" = 8074.314724951111; "R" = 72668.83252455998; "S" = 0.1620203638381736; "T" = 0.1620203638381736 }
    9 = @{ "E" = 3; "G" = 110.1980973333333; "H" = 330.594292; "I" = 0.2179577200213544; "J" = 0.2179577200213544; "K" = 3; "M" = 2.513102666666667; "N" = 7.539308; "O" = 0.02549623188763971; "P" = 0.02549623188763971; "Q" = 276.9391322699929; "R" = 2492.452190429936; "S" = 0.005557100571365705; "T" = 0.005557100571365705 }
    10 = @{ "E" = 3; "G" = 351.7202226666666; "H" = 1055.160668; "I" = 0.6956575446665283; "J" = 0.6956575446665284; "K" = 3; "M" = 11.84297166666667; "N" = 35.528915; "O" = 0.1201507426883529; "P" = 0.1201507426883529; "Q" = 4165.412631635024; "R" = 37488.71368471521; "S" = 0.08358377064843941; "T" = 0.08358377064843944 }
    11 = @{ "E" = 3; "G" = 351.7202226666666; "H" = 1055.160668; "I" = 0.6956575446665283; "J" = 0.6956575446665284; "K" = 3; "M" = 10.940628; "N" = 32.821884; "O" = 0.1109961770302011; "P" = 0.1109961770302011; "Q" = 3848.040116273167; "R" = 34632.36104645851; "S" = 0.077215327980201; "T" = 0.07721532798020103 }
    12 = @{ "E" = 3; "G" = 351.7202226666666; "H" = 1055.160668; "I" = 0.6956575446665283; "J" = 0.6956575446665284; "K" = 3; "M" = 73.27090866666667; "N" = 219.812726; "O" = 0.7433568483938062; "P" = 0.7433568483938062; "Q" = 25770.860311229; "R" = 231937.742801061; "S" = 0.517121799964684; "T" = 0.517121799964684 }
    13 = @{ "E" = 3; "G" = 351.7202226666666; "H" = 1055.160668; "I" = 0.6956575446665283; "J" = 0.6956575446665284; "K" = 3; "M" = 2.513102666666667; "N" = 7.539308; "O" = 0.02549623188763971; "P" = 0.02549623188763971; "Q" = 883.9090295041938; "R" = 7955.181265537744; "S" = 0.01773664607320389; "T" = 0.01773664607320389 }
    14 = @{ "E" = 3; "G" = 19.17031833333333; "H" = 57.510955; "I" = 0.03791643392333802; "J" = 0.03791643392333802; "K" = 3; "M" = 11.84297166666667; "N" = 35.528915; "O" = 0.1201507426883529; "P" = 0.1201507426883529; "Q" = 227.0335368626472; "R" = 2043.301831763825; "S" = 0.004555687695982922; "T" = 0.004555687695982922 }
    15 = @{ "E" = 3; "G" = 19.17031833333333; "H" = 57.510955; "I" = 0.03791643392333802; "J" = 0.03791643392333802; "K" = 3; "M" = 10.940628; "N" = 32.821884; "O" = 0.1109961770302011; "P" = 0.1109961770302011; "Q" = 209.73532152658; "R" = 1887.61789373922; "S" = 0.004208579212108749; "T" = 0.004208579212108749 }
    16 = @{ "E" = 3; "G" = 19.17031833333333; "H" = 57.510955; "I" = 0.03791643392333802; "J" = 0.03791643392333802; "K" = 3; "M" = 73.27090866666667; "N" = 219.812726; "O" = 0.7433568483938062; "P" = 0.7433568483938062; "Q" = 1404.626643712592; "R" = 12641.63979341333; "S" = 0.02818544082358455; "T" = 0.02818544082358455 }
    17 = @{ "E" = 3; "G" = 19.17031833333333; "H" = 57.510955; "I" = 0.03791643392333802; "J" = 0.03791643392333802; "K" = 3; "M" = 2.513102666666667; "N" = 7.539308; "O" = 0.02549623188763971; "P" = 0.02549623188763971; "Q" = 48.17697812434889; "R" = 433.59280311914; "S" = 0.0009667261916617947; "T" = 0.0009667261916617947 }
}

foreach ($r in $rowData.Keys) {
    foreach ($col in $rowData[$r].Keys) {
        $ws.Range("$col$r").Value = $rowData[$r][$col]
    }
}